$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new first column for the "Time Stamp" values; this shifts the
# existing Temperature/Humidity/Pressure/Wind Speed columns from A:D to B:E.
$ws.Columns.Item(1).Insert()

# Header row
$ws.Cells.Item(1, 1).Value = "Time Stamp"

# Remove the old blank separator rows (3 and 5) and the now-redundant
# rows 6 and 7 so only two data rows remain.
$ws.Rows.Item(3).Delete()
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(4).Delete()

# Row 2 - first reading (no wind speed reported)
$ws.Cells.Item(2, 1).Value = "2017.05.19 16.07.33"
$ws.Cells.Item(2, 2).Value = 301.0899963378906
$ws.Cells.Item(2, 3).Value = 30.0
$ws.Cells.Item(2, 4).Value = 1012.0
$ws.Cells.Item(2, 5).ClearContents()

# Row 3 - second reading
$ws.Cells.Item(3, 1).Value = "2017.05.19 16.08.17"
$ws.Cells.Item(3, 2).Value = 301.0899963378906
$ws.Cells.Item(3, 3).Value = 30.0
$ws.Cells.Item(3, 4).Value = 1012.0
$ws.Cells.Item(3, 5).Value = 4.599999904632568
